$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column V: header + sample phone number value
$ws.Range("V1").Value = "pp_phone_no_i_c"
$ws.Range("V2").Value = "+48603499023"

# Adjust view / selection to match target state
$ws.Range("V3").Select()
